$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# ("1.003", "8.448", etc.) are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "27.169.48"
$ws.Cells.Item(2, 5).Value = "  -1.09%  "

$ws.Cells.Item(3, 4).Value = "1.784.04"
$ws.Cells.Item(3, 5).Value = "  -1.79%  "

$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "

$ws.Cells.Item(5, 4).Value = "337.49"
$ws.Cells.Item(5, 5).Value = "  -1.80%  "

$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.15%  "

$ws.Cells.Item(7, 4).Value = "0.3861"
$ws.Cells.Item(7, 5).Value = "  +0.90%  "

$ws.Cells.Item(8, 4).Value = "0.3433"
$ws.Cells.Item(8, 5).Value = "  -2.33%  "

$ws.Cells.Item(9, 4).Value = "47.90"
$ws.Cells.Item(9, 5).Value = "  -2.08%  "

$ws.Cells.Item(10, 4).Value = "1.191"
$ws.Cells.Item(10, 5).Value = "  -3.44%  "

$ws.Cells.Item(11, 4).Value = "0.07456"
$ws.Cells.Item(11, 5).Value = "  -4.68%  "

$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  -0.06%  "

$ws.Cells.Item(13, 4).Value = "21.69"
$ws.Cells.Item(13, 5).Value = "  -2.52%  "

$ws.Cells.Item(14, 4).Value = "6.444"
$ws.Cells.Item(14, 5).Value = "  -2.40%  "

$ws.Cells.Item(15, 4).Value = "1.782.06"
$ws.Cells.Item(15, 5).Value = "  -2.03%  "

$ws.Cells.Item(16, 4).Value = "7.127"
$ws.Cells.Item(16, 5).Value = "  -1.48%  "

$ws.Cells.Item(17, 4).Value = "0.00001092"
$ws.Cells.Item(17, 5).Value = "  -2.57%  "

$ws.Cells.Item(18, 4).Value = "0.06661"
$ws.Cells.Item(18, 5).Value = "  -0.99%  "

$ws.Cells.Item(19, 4).Value = "83.38"
$ws.Cells.Item(19, 5).Value = "  -3.46%  "

$ws.Cells.Item(20, 5).Value = "  +0.11%  "

$ws.Cells.Item(21, 4).Value = "17.56"
$ws.Cells.Item(21, 5).Value = "  -0.58%  "

$ws.Cells.Item(22, 4).Value = "6.516"
$ws.Cells.Item(22, 5).Value = "  -0.95%  "

$ws.Cells.Item(23, 4).Value = "27.162.22"
$ws.Cells.Item(23, 5).Value = "  -1.16%  "

$ws.Cells.Item(24, 4).Value = "12.35"
$ws.Cells.Item(24, 5).Value = "  -6.25%  "

$ws.Cells.Item(25, 4).Value = "2.366"
$ws.Cells.Item(25, 5).Value = "  -3.83%  "

$ws.Cells.Item(26, 4).Value = "21.15"
$ws.Cells.Item(26, 5).Value = "  -4.67%  "

$ws.Cells.Item(27, 4).Value = "2.497"
$ws.Cells.Item(27, 5).Value = "  -6.76%  "

$ws.Cells.Item(28, 4).Value = "1.446"
$ws.Cells.Item(28, 5).Value = "  -1.60%  "

$ws.Cells.Item(29, 4).Value = "156.59"
$ws.Cells.Item(29, 5).Value = "  +1.65%  "

$ws.Cells.Item(30, 4).Value = "1.983.13"
$ws.Cells.Item(30, 5).Value = "  -2.01%  "

$ws.Cells.Item(31, 4).Value = "134.24"
$ws.Cells.Item(31, 5).Value = "  -1.57%  "

$ws.Cells.Item(32, 4).Value = "3.972"
$ws.Cells.Item(32, 5).Value = "  -2.22%  "

$ws.Cells.Item(33, 4).Value = "5.996"

$ws.Cells.Item(34, 4).Value = "0.08694"
$ws.Cells.Item(34, 5).Value = "  -1.24%  "

$ws.Cells.Item(35, 4).Value = "12.97"
$ws.Cells.Item(35, 5).Value = "  -6.58%  "

$ws.Cells.Item(36, 4).Value = "1.621"
$ws.Cells.Item(36, 5).Value = "  -4.07%  "

$ws.Cells.Item(37, 4).Value = "5.403"
$ws.Cells.Item(37, 5).Value = "  -4.08%  "

$ws.Cells.Item(38, 4).Value = "0.6823"

$ws.Cells.Item(39, 4).Value = "0.06340"
$ws.Cells.Item(39, 5).Value = "  -2.22%  "

$ws.Cells.Item(40, 4).Value = "0.02346"
$ws.Cells.Item(40, 5).Value = "  -2.59%  "

$ws.Cells.Item(41, 4).Value = "0.2191"
$ws.Cells.Item(41, 5).Value = "  -3.27%  "

$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "8.448"
$ws.Cells.Item(42, 5).Value = "  -5.76%  "

$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).Value = "1.236"
$ws.Cells.Item(43, 5).Value = "  -5.13%  "

$ws.Cells.Item(44, 4).Value = "14.24"
$ws.Cells.Item(44, 5).Value = "  -3.89%  "

$ws.Cells.Item(45, 4).Value = "1.000"
$ws.Cells.Item(45, 5).Value = "  +0.08%  "

$ws.Cells.Item(46, 4).Value = "0.6408"
$ws.Cells.Item(46, 5).Value = "  -2.53%  "

$ws.Cells.Item(47, 4).Value = "3.859"
$ws.Cells.Item(47, 5).Value = "  -2.57%  "

$ws.Cells.Item(48, 4).Value = "2.170"
$ws.Cells.Item(48, 5).Value = "  -0.59%  "

$ws.Cells.Item(49, 4).Value = "131.59"
$ws.Cells.Item(49, 5).Value = "  -0.88%  "

$ws.Cells.Item(50, 4).Value = "0.07112"
$ws.Cells.Item(50, 5).Value = "  -2.90%  "

$ws.Cells.Item(51, 4).Value = "79.38"
$ws.Cells.Item(51, 5).Value = "  -1.46%  "
